$wb = $excel.ActiveWorkbook

# --- 1. Remove the stray empty cell B2 on "ODI Batting" ---
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B2").Value = $null

# --- 2. Add the new "ODI Batting Extra" sheet after the last existing sheet ---
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$headerStyleSource = $odiBowling.Range("A1:F1")

$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "ODI Batting Extra"

# Reuse the same bold/centered header style used by the other sheets.
$headerStyleSource.Copy($ws.Range("A1:F1"))

# Header row
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Data row (A2 must stay text, matching MATCH_CODE's representation elsewhere)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "4460"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = $null
$ws.Range("D2").Value = $null
$ws.Range("E2").Value = $null
$ws.Range("F2").Value = "NO"
